$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Min / Max / Average of Quantity (E) and Order Value (F)
$ws.Range("I4").Formula = '=MIN($E:$E)'
$ws.Range("J4").Formula = '=MIN($F$2:$F$30)'

$ws.Range("I5").Formula = '=MAX($E:$E)'
$ws.Range("J5").Formula = '=MAX($F$2:$F$30)'

$ws.Range("I6").Formula = '=AVERAGE($E:$E)'
$ws.Range("J6").Formula = '=AVERAGE($F$2:$F$30)'

# Average / Min / Max Order Value per Product Type
$ws.Range("I13").Formula = '=AVERAGEIF($C$2:$C$30, $H13, $F$2:$F$30)'
$ws.Range("J13").Formula = '=MINIFS($F$2:$F$30,$C$2:$C$30,H13)'
$ws.Range("K13").Formula = '=MAXIFS($F$2:$F$30,$C$2:$C$30,H13)'

$ws.Range("I14").NumberFormat = "0"
$ws.Range("I14").Formula = '=AVERAGEIF($C$2:$C$30, H14, $F$2:$F$30)'
$ws.Range("J14").NumberFormat = "0"
$ws.Range("J14").Formula = '=MINIFS($F$2:$F$30,$C$2:$C$30,H14)'
$ws.Range("K14").NumberFormat = "0"
$ws.Range("K14").Formula = '=MAXIFS($F$2:$F$30,$C$2:$C$30,H14)'

$ws.Range("I15").NumberFormat = "0"
$ws.Range("I15").Formula = '=AVERAGEIF($C$2:$C$30, $H15, $F$2:$F$30)'
$ws.Range("J15").Formula = '=MINIFS($F$2:$F$30,$C$2:$C$30,H15)'
$ws.Range("K15").Formula = '=MAXIFS($F$2:$F$30,$C$2:$C$30,H15)'

# Update the active selection to match the final state of the workbook
$ws.Range("N19").Select()
